$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2,1).Value = "ECs"
$ws.Cells.Item(2,2).Value = "Vegfc"
$ws.Cells.Item(2,3).Value = "Vipr2"
$ws.Cells.Item(2,4).Value = "ECs"
$ws.Cells.Item(2,5).Value = 3.0
$ws.Cells.Item(2,6).Value = 1.0
$ws.Cells.Item(2,7).Value = 2.404594333333333
$ws.Cells.Item(2,8).Value = 7.213783
$ws.Cells.Item(2,9).Value = 0.3565065439253589
$ws.Cells.Item(2,10).Value = 0.3565065439253589
$ws.Cells.Item(2,11).Value = 1.0
$ws.Cells.Item(2,12).Value = 0.3333333333333333
$ws.Cells.Item(2,13).Value = 0.01962266666666667
$ws.Cells.Item(2,14).Value = 0.058868
$ws.Cells.Item(2,15).Value = 0.02207703711370904
$ws.Cells.Item(2,16).Value = 0.02207703711370904
$ws.Cells.Item(2,17).Value = 0.04718455307155556
$ws.Cells.Item(2,18).Value = 0.424660977644
$ws.Cells.Item(2,19).Value = 0.007870608201520291
$ws.Cells.Item(2,20).Value = 0.00787060820152029

# Row 3
$ws.Cells.Item(3,1).Value = "ECs"
$ws.Cells.Item(3,2).Value = "Vegfc"
$ws.Cells.Item(3,3).Value = "Vipr2"
$ws.Cells.Item(3,4).Value = "FAPs"
$ws.Cells.Item(3,5).Value = 3.0
$ws.Cells.Item(3,6).Value = 1.0
$ws.Cells.Item(3,7).Value = 2.404594333333333
$ws.Cells.Item(3,8).Value = 7.213783
$ws.Cells.Item(3,9).Value = 0.3565065439253589
$ws.Cells.Item(3,10).Value = 0.3565065439253589
$ws.Cells.Item(3,11).Value = 3.0
$ws.Cells.Item(3,12).Value = 1.0
$ws.Cells.Item(3,13).Value = 0.6658376666666667
$ws.Cells.Item(3,14).Value = 1.997513
$ws.Cells.Item(3,15).Value = 0.7491195324474467
$ws.Cells.Item(3,16).Value = 0.7491195324474467
$ws.Cells.Item(3,17).Value = 1.601069480186555
$ws.Cells.Item(3,18).Value = 14.409625321679
$ws.Cells.Item(3,19).Value = 0.26706601549982
$ws.Cells.Item(3,20).Value = 0.26706601549982

# Row 4
$ws.Cells.Item(4,1).Value = "ECs"
$ws.Cells.Item(4,2).Value = "Vegfc"
$ws.Cells.Item(4,3).Value = "Vipr2"
$ws.Cells.Item(4,4).Value = "Inflammatory-Mac"
$ws.Cells.Item(4,5).Value = 3.0
$ws.Cells.Item(4,6).Value = 1.0
$ws.Cells.Item(4,7).Value = 2.404594333333333
$ws.Cells.Item(4,8).Value = 7.213783
$ws.Cells.Item(4,9).Value = 0.3565065439253589
$ws.Cells.Item(4,10).Value = 0.3565065439253589
$ws.Cells.Item(4,11).Value = 1.0
$ws.Cells.Item(4,12).Value = 0.3333333333333333
$ws.Cells.Item(4,13).Value = 0.1131433333333333
$ws.Cells.Item(4,14).Value = 0.33943
$ws.Cells.Item(4,15).Value = 0.1272951129222372
$ws.Cells.Item(4,16).Value = 0.1272951129222372
$ws.Cells.Item(4,17).Value = 0.2720638181877778
$ws.Cells.Item(4,18).Value = 2.44857436369
$ws.Cells.Item(4,19).Value = 0.04538154076649507
$ws.Cells.Item(4,20).Value = 0.04538154076649508

# Row 5
$ws.Cells.Item(5,1).Value = "ECs"
$ws.Cells.Item(5,2).Value = "Vegfc"
$ws.Cells.Item(5,3).Value = "Vipr2"
$ws.Cells.Item(5,4).Value = "MuSCs"
$ws.Cells.Item(5,5).Value = 3.0
$ws.Cells.Item(5,6).Value = 1.0
$ws.Cells.Item(5,7).Value = 2.404594333333333
$ws.Cells.Item(5,8).Value = 7.213783
$ws.Cells.Item(5,9).Value = 0.3565065439253589
$ws.Cells.Item(5,10).Value = 0.3565065439253589
$ws.Cells.Item(5,11).Value = 2.0
$ws.Cells.Item(5,12).Value = 0.6666666666666666
$ws.Cells.Item(5,13).Value = 0.09022333333333332
$ws.Cells.Item(5,14).Value = 0.27067
$ws.Cells.Item(5,15).Value = 0.1015083175166071
$ws.Cells.Item(5,16).Value = 0.1015083175166071
$ws.Cells.Item(5,17).Value = 0.2169505160677777
$ws.Cells.Item(5,18).Value = 1.95255464461
$ws.Cells.Item(5,19).Value = 0.03618837945752355
$ws.Cells.Item(5,20).Value = 0.03618837945752355

# Row 6
$ws.Cells.Item(6,1).Value = "FAPs"
$ws.Cells.Item(6,2).Value = "Vegfc"
$ws.Cells.Item(6,3).Value = "Vipr2"
$ws.Cells.Item(6,4).Value = "ECs"
$ws.Cells.Item(6,5).Value = 3.0
$ws.Cells.Item(6,6).Value = 1.0
$ws.Cells.Item(6,7).Value = 3.178631333333334
$ws.Cells.Item(6,8).Value = 9.535894
$ws.Cells.Item(6,9).Value = 0.471265716362492
$ws.Cells.Item(6,10).Value = 0.4712657163624919
$ws.Cells.Item(6,11).Value = 1.0
$ws.Cells.Item(6,12).Value = 0.3333333333333333
$ws.Cells.Item(6,13).Value = 0.01962266666666667
$ws.Cells.Item(6,14).Value = 0.058868
$ws.Cells.Item(6,15).Value = 0.02207703711370904
$ws.Cells.Item(6,16).Value = 0.02207703711370904
$ws.Cells.Item(6,17).Value = 0.06237322311022223
$ws.Cells.Item(6,18).Value = 0.561359007992
$ws.Cells.Item(6,19).Value = 0.01040415071055341
$ws.Cells.Item(6,20).Value = 0.01040415071055341

# Row 7
$ws.Cells.Item(7,1).Value = "FAPs"
$ws.Cells.Item(7,2).Value = "Vegfc"
$ws.Cells.Item(7,3).Value = "Vipr2"
$ws.Cells.Item(7,4).Value = "FAPs"
$ws.Cells.Item(7,5).Value = 3.0
$ws.Cells.Item(7,6).Value = 1.0
$ws.Cells.Item(7,7).Value = 3.178631333333334
$ws.Cells.Item(7,8).Value = 9.535894
$ws.Cells.Item(7,9).Value = 0.471265716362492
$ws.Cells.Item(7,10).Value = 0.4712657163624919
$ws.Cells.Item(7,11).Value = 3.0
$ws.Cells.Item(7,12).Value = 1.0
$ws.Cells.Item(7,13).Value = 0.6658376666666667
$ws.Cells.Item(7,14).Value = 1.997513
$ws.Cells.Item(7,15).Value = 0.7491195324474467
$ws.Cells.Item(7,16).Value = 0.7491195324474467
$ws.Cells.Item(7,17).Value = 2.116452470180223
$ws.Cells.Item(7,18).Value = 19.048072231622
$ws.Cells.Item(7,19).Value = 0.353034353099981
$ws.Cells.Item(7,20).Value = 0.353034353099981

# Row 8
$ws.Cells.Item(8,1).Value = "FAPs"
$ws.Cells.Item(8,2).Value = "Vegfc"
$ws.Cells.Item(8,3).Value = "Vipr2"
$ws.Cells.Item(8,4).Value = "Inflammatory-Mac"
$ws.Cells.Item(8,5).Value = 3.0
$ws.Cells.Item(8,6).Value = 1.0
$ws.Cells.Item(8,7).Value = 3.178631333333334
$ws.Cells.Item(8,8).Value = 9.535894
$ws.Cells.Item(8,9).Value = 0.471265716362492
$ws.Cells.Item(8,10).Value = 0.4712657163624919
$ws.Cells.Item(8,11).Value = 1.0
$ws.Cells.Item(8,12).Value = 0.3333333333333333
$ws.Cells.Item(8,13).Value = 0.1131433333333333
$ws.Cells.Item(8,14).Value = 0.33943
$ws.Cells.Item(8,15).Value = 0.1272951129222372
$ws.Cells.Item(8,16).Value = 0.1272951129222372
$ws.Cells.Item(8,17).Value = 0.3596409444911111
$ws.Cells.Item(8,18).Value = 3.23676850042
$ws.Cells.Item(8,19).Value = 0.05998982258074242
$ws.Cells.Item(8,20).Value = 0.05998982258074242

# Row 9
$ws.Cells.Item(9,1).Value = "FAPs"
$ws.Cells.Item(9,2).Value = "Vegfc"
$ws.Cells.Item(9,3).Value = "Vipr2"
$ws.Cells.Item(9,4).Value = "MuSCs"
$ws.Cells.Item(9,5).Value = 3.0
$ws.Cells.Item(9,6).Value = 1.0
$ws.Cells.Item(9,7).Value = 3.178631333333334
$ws.Cells.Item(9,8).Value = 9.535894
$ws.Cells.Item(9,9).Value = 0.471265716362492
$ws.Cells.Item(9,10).Value = 0.4712657163624919
$ws.Cells.Item(9,11).Value = 2.0
$ws.Cells.Item(9,12).Value = 0.6666666666666666
$ws.Cells.Item(9,13).Value = 0.09022333333333332
$ws.Cells.Item(9,14).Value = 0.27067
$ws.Cells.Item(9,15).Value = 0.1015083175166071
$ws.Cells.Item(9,16).Value = 0.1015083175166071
$ws.Cells.Item(9,17).Value = 0.2867867143311111
$ws.Cells.Item(9,18).Value = 2.58108042898
$ws.Cells.Item(9,19).Value = 0.04783738997121513
$ws.Cells.Item(9,20).Value = 0.04783738997121512

# Row 10
$ws.Cells.Item(10,1).Value = "MuSCs"
$ws.Cells.Item(10,2).Value = "Vegfc"
$ws.Cells.Item(10,3).Value = "Vipr2"
$ws.Cells.Item(10,4).Value = "ECs"
$ws.Cells.Item(10,5).Value = 3.0
$ws.Cells.Item(10,6).Value = 1.0
$ws.Cells.Item(10,7).Value = 1.125649
$ws.Cells.Item(10,8).Value = 3.376947
$ws.Cells.Item(10,9).Value = 0.166889370527102
$ws.Cells.Item(10,10).Value = 0.1668893705271019
$ws.Cells.Item(10,11).Value = 1.0
$ws.Cells.Item(10,12).Value = 0.3333333333333333
$ws.Cells.Item(10,13).Value = 0.01962266666666667
$ws.Cells.Item(10,14).Value = 0.058868
$ws.Cells.Item(10,15).Value = 0.02207703711370904
$ws.Cells.Item(10,16).Value = 0.02207703711370904
$ws.Cells.Item(10,17).Value = 0.02208823511066667
$ws.Cells.Item(10,18).Value = 0.198794115996
$ws.Cells.Item(10,19).Value = 0.00368442282701037
$ws.Cells.Item(10,20).Value = 0.003684422827010369

# Row 11
$ws.Cells.Item(11,1).Value = "MuSCs"
$ws.Cells.Item(11,2).Value = "Vegfc"
$ws.Cells.Item(11,3).Value = "Vipr2"
$ws.Cells.Item(11,4).Value = "FAPs"
$ws.Cells.Item(11,5).Value = 3.0
$ws.Cells.Item(11,6).Value = 1.0
$ws.Cells.Item(11,7).Value = 1.125649
$ws.Cells.Item(11,8).Value = 3.376947
$ws.Cells.Item(11,9).Value = 0.166889370527102
$ws.Cells.Item(11,10).Value = 0.1668893705271019
$ws.Cells.Item(11,11).Value = 3.0
$ws.Cells.Item(11,12).Value = 1.0
$ws.Cells.Item(11,13).Value = 0.6658376666666667
$ws.Cells.Item(11,14).Value = 1.997513
$ws.Cells.Item(11,15).Value = 0.7491195324474467
$ws.Cells.Item(11,16).Value = 0.7491195324474467
$ws.Cells.Item(11,17).Value = 0.7494995036456668
$ws.Cells.Item(11,18).Value = 6.745495532811001
$ws.Cells.Item(11,19).Value = 0.1250200872197113
$ws.Cells.Item(11,20).Value = 0.1250200872197113

# Row 12
$ws.Cells.Item(12,1).Value = "MuSCs"
$ws.Cells.Item(12,2).Value = "Vegfc"
$ws.Cells.Item(12,3).Value = "Vipr2"
$ws.Cells.Item(12,4).Value = "Inflammatory-Mac"
$ws.Cells.Item(12,5).Value = 3.0
$ws.Cells.Item(12,6).Value = 1.0
$ws.Cells.Item(12,7).Value = 1.125649
$ws.Cells.Item(12,8).Value = 3.376947
$ws.Cells.Item(12,9).Value = 0.166889370527102
$ws.Cells.Item(12,10).Value = 0.1668893705271019
$ws.Cells.Item(12,11).Value = 1.0
$ws.Cells.Item(12,12).Value = 0.3333333333333333
$ws.Cells.Item(12,13).Value = 0.1131433333333333
$ws.Cells.Item(12,14).Value = 0.33943
$ws.Cells.Item(12,15).Value = 0.1272951129222372
$ws.Cells.Item(12,16).Value = 0.1272951129222372
$ws.Cells.Item(12,17).Value = 0.1273596800233333
$ws.Cells.Item(12,18).Value = 1.14623712021
$ws.Cells.Item(12,19).Value = 0.02124420126676853
$ws.Cells.Item(12,20).Value = 0.02124420126676853

# Row 13
$ws.Cells.Item(13,1).Value = "MuSCs"
$ws.Cells.Item(13,2).Value = "Vegfc"
$ws.Cells.Item(13,3).Value = "Vipr2"
$ws.Cells.Item(13,4).Value = "MuSCs"
$ws.Cells.Item(13,5).Value = 3.0
$ws.Cells.Item(13,6).Value = 1.0
$ws.Cells.Item(13,7).Value = 1.125649
$ws.Cells.Item(13,8).Value = 3.376947
$ws.Cells.Item(13,9).Value = 0.166889370527102
$ws.Cells.Item(13,10).Value = 0.1668893705271019
$ws.Cells.Item(13,11).Value = 2.0
$ws.Cells.Item(13,12).Value = 0.6666666666666666
$ws.Cells.Item(13,13).Value = 0.09022333333333332
$ws.Cells.Item(13,14).Value = 0.27067
$ws.Cells.Item(13,15).Value = 0.1015083175166071
$ws.Cells.Item(13,16).Value = 0.1015083175166071
$ws.Cells.Item(13,17).Value = 0.1015598049433333
$ws.Cells.Item(13,18).Value = 0.9140382444899999
$ws.Cells.Item(13,19).Value = 0.01694065921361175
$ws.Cells.Item(13,20).Value = 0.01694065921361175

# Row 14
$ws.Cells.Item(14,1).Value = "Neutrophils"
$ws.Cells.Item(14,2).Value = "Vegfc"
$ws.Cells.Item(14,3).Value = "Vipr2"
$ws.Cells.Item(14,4).Value = "ECs"
$ws.Cells.Item(14,5).Value = 1.0
$ws.Cells.Item(14,6).Value = 0.3333333333333333
$ws.Cells.Item(14,7).Value = 0.03600666666666667
$ws.Cells.Item(14,8).Value = 0.10802
$ws.Cells.Item(14,9).Value = 0.005338369185047189
$ws.Cells.Item(14,10).Value = 0.005338369185047188
$ws.Cells.Item(14,11).Value = 1.0
$ws.Cells.Item(14,12).Value = 0.3333333333333333
$ws.Cells.Item(14,13).Value = 0.01962266666666667
$ws.Cells.Item(14,14).Value = 0.058868
$ws.Cells.Item(14,15).Value = 0.02207703711370904
$ws.Cells.Item(14,16).Value = 0.02207703711370904
$ws.Cells.Item(14,17).Value = 0.0007065468177777777
$ws.Cells.Item(14,18).Value = 0.00635892136
$ws.Cells.Item(14,19).Value = 0.0001178553746249675
$ws.Cells.Item(14,20).Value = 0.0001178553746249675

# Row 15
$ws.Cells.Item(15,1).Value = "Neutrophils"
$ws.Cells.Item(15,2).Value = "Vegfc"
$ws.Cells.Item(15,3).Value = "Vipr2"
$ws.Cells.Item(15,4).Value = "FAPs"
$ws.Cells.Item(15,5).Value = 1.0
$ws.Cells.Item(15,6).Value = 0.3333333333333333
$ws.Cells.Item(15,7).Value = 0.03600666666666667
$ws.Cells.Item(15,8).Value = 0.10802
$ws.Cells.Item(15,9).Value = 0.005338369185047189
$ws.Cells.Item(15,10).Value = 0.005338369185047188
$ws.Cells.Item(15,11).Value = 3.0
$ws.Cells.Item(15,12).Value = 1.0
$ws.Cells.Item(15,13).Value = 0.6658376666666667
$ws.Cells.Item(15,14).Value = 1.997513
$ws.Cells.Item(15,15).Value = 0.7491195324474467
$ws.Cells.Item(15,16).Value = 0.7491195324474467
$ws.Cells.Item(15,17).Value = 0.02397459491777778
$ws.Cells.Item(15,18).Value = 0.21577135426
$ws.Cells.Item(15,19).Value = 0.003999076627934407
$ws.Cells.Item(15,20).Value = 0.003999076627934407

# Row 16
$ws.Cells.Item(16,1).Value = "Neutrophils"
$ws.Cells.Item(16,2).Value = "Vegfc"
$ws.Cells.Item(16,3).Value = "Vipr2"
$ws.Cells.Item(16,4).Value = "Inflammatory-Mac"
$ws.Cells.Item(16,5).Value = 1.0
$ws.Cells.Item(16,6).Value = 0.3333333333333333
$ws.Cells.Item(16,7).Value = 0.03600666666666667
$ws.Cells.Item(16,8).Value = 0.10802
$ws.Cells.Item(16,9).Value = 0.005338369185047189
$ws.Cells.Item(16,10).Value = 0.005338369185047188
$ws.Cells.Item(16,11).Value = 1.0
$ws.Cells.Item(16,12).Value = 0.3333333333333333
$ws.Cells.Item(16,13).Value = 0.1131433333333333
$ws.Cells.Item(16,14).Value = 0.33943
$ws.Cells.Item(16,15).Value = 0.1272951129222372
$ws.Cells.Item(16,16).Value = 0.1272951129222372
$ws.Cells.Item(16,17).Value = 0.004073914288888889
$ws.Cells.Item(16,18).Value = 0.0366652286
$ws.Cells.Item(16,19).Value = 0.0006795483082311733
$ws.Cells.Item(16,20).Value = 0.0006795483082311733

# Row 17
$ws.Cells.Item(17,1).Value = "Neutrophils"
$ws.Cells.Item(17,2).Value = "Vegfc"
$ws.Cells.Item(17,3).Value = "Vipr2"
$ws.Cells.Item(17,4).Value = "MuSCs"
$ws.Cells.Item(17,5).Value = 1.0
$ws.Cells.Item(17,6).Value = 0.3333333333333333
$ws.Cells.Item(17,7).Value = 0.03600666666666667
$ws.Cells.Item(17,8).Value = 0.10802
$ws.Cells.Item(17,9).Value = 0.005338369185047189
$ws.Cells.Item(17,10).Value = 0.005338369185047188
$ws.Cells.Item(17,11).Value = 2.0
$ws.Cells.Item(17,12).Value = 0.6666666666666666
$ws.Cells.Item(17,13).Value = 0.09022333333333332
$ws.Cells.Item(17,14).Value = 0.27067
$ws.Cells.Item(17,15).Value = 0.1015083175166071
$ws.Cells.Item(17,16).Value = 0.1015083175166071
$ws.Cells.Item(17,17).Value = 0.003248641488888889
$ws.Cells.Item(17,18).Value = 0.0292377734
$ws.Cells.Item(17,19).Value = 0.0005418888742566411
$ws.Cells.Item(17,20).Value = 0.000541888874256641

